$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update buyer-name (F) and buyer-email (E) for rows 15-17: Smith Johnson -> Dave Harley Petterson
# (set the name before the email so new shared strings are appended name-first,
# matching the order they were added upstream)
$ws.Range("F15").Value = "Dave Harley Petterson"
$ws.Range("E15").Value = "daveharleypetterson02@gmail.com"

$ws.Range("F16").Value = "Dave Harley Petterson"
$ws.Range("E16").Value = "daveharleypetterson02@gmail.com"

$ws.Range("F17").Value = "Dave Harley Petterson"
$ws.Range("E17").Value = "daveharleypetterson02@gmail.com"

# Hyperlink for E15 (single cell, no display text - mirrors E9/E12/E27 pattern)
$hE15 = $ws.Hyperlinks.Add($ws.Range("E15"), "mailto:daveharleypetterson02@gmail.com")

# Hyperlink for E16:E17 as one merged range (mirrors E10:E11 / E13:E14 / E28:E29 pattern).
# Add per-cell first so both cells pick up the hyperlink style, then delete those
# and add a single combined hyperlink so the saved <hyperlinks> entry covers the
# whole E16:E17 range (matching how the existing merged hyperlinks look).
$hE16 = $ws.Hyperlinks.Add($ws.Range("E16"), "mailto:daveharleypetterson02@gmail.com")
$hE17 = $ws.Hyperlinks.Add($ws.Range("E17"), "mailto:daveharleypetterson02@gmail.com")
$hE16.Delete()
$hE17.Delete()
$ws.Hyperlinks.Add($ws.Range("E16:E17"), "mailto:daveharleypetterson02@gmail.com", "", "", "daveharleypetterson02@gmail.com")

# Widen column E to fit the longer email address and drop the old bestFit autosizing
$ws.Columns.Item(5).ColumnWidth = 46.8

# Move the active selection from G29 to G9
$ws.Range("G9").Select() | Out-Null
